$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "RR"
$ws.Range("A3").Value = "MT"
$ws.Range("A4").Value = "PI"
$ws.Range("A5").Value = "TO"
$ws.Range("A6").Value = "AC"
$ws.Range("A7").Value = "PB"
$ws.Range("A8").Value = "SE"
$ws.Range("A9").Value = "BR"
$ws.Range("A10").Value = "NE"

$ws.Range("D2").Value = "RR"
$ws.Range("D3").Value = "MT"
$ws.Range("D4").Value = "TO"
$ws.Range("D5").Value = "PI"
$ws.Range("D6").Value = "MS"
$ws.Range("D7").Value = "MA"
$ws.Range("D8").Value = "SE"
$ws.Range("D9").Value = "BR"
$ws.Range("D10").Value = "NE"
